$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to reflect repulled data / recalculated means
$ws.Range("F2").Value = -6
$ws.Range("F4").Value = 8
$ws.Range("F5").Value = -11
$ws.Range("F7").Value = -2
